# Update rows 2-10 of Sheet1 with recomputed NATMI ligand-receptor statistics
# (per Dr Hou's advice), expanding the target-cluster cross join to include sCs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("ECs", "Il15", "Il2ra", "ECs", 3, 1, 16.58023833333333, 49.740715, 0.63541025828417, 0.63541025828417, 2, 0.6666666666666666, 0.3213843333333333, 0.964153, 0.1066545648432073, 0.1066545648432073, 5.328628843266111, 47.957659589395, 0.06776940459420808, 0.06776940459420808)
$row3 = @("ECs", "Il15", "Il2ra", "FAPs", 3, 1, 16.58023833333333, 49.740715, 0.63541025828417, 0.63541025828417, 3, 1, 1.174648, 3.523944, 0.3898185390200842, 0.3898185390200842, 19.47594379777334, 175.28349417996, 0.2476946985627095, 0.2476946985627095)
$row4 = @("ECs", "Il15", "Il2ra", "sCs", 3, 1, 16.58023833333333, 49.740715, 0.63541025828417, 0.63541025828417, 3, 1, 1.517287666666667, 4.551863, 0.5035268961367085, 0.5035268961367085, 25.15699113356056, 226.412920202045, 0.3199461551272524, 0.3199461551272524)
$row5 = @("FAPs", "Il15", "Il2ra", "ECs", 3, 1, 8.574149, 25.722447, 0.3285901035393414, 0.3285901035393414, 2, 0.6666666666666666, 0.3213843333333333, 0.964153, 0.1066545648432073, 0.1066545648432073, 2.755597160265667, 24.800374442391, 0.03504563450477288, 0.03504563450477288)
$row6 = @("FAPs", "Il15", "Il2ra", "FAPs", 3, 1, 8.574149, 25.722447, 0.3285901035393414, 0.3285901035393414, 3, 1, 1.174648, 3.523944, 0.3898185390200842, 0.3898185390200842, 10.071606974552, 90.644462770968, 0.1280905140981643, 0.1280905140981643)
$row7 = @("FAPs", "Il15", "Il2ra", "sCs", 3, 1, 8.574149, 25.722447, 0.3285901035393414, 0.3285901035393414, 3, 1, 1.517287666666667, 4.551863, 0.5035268961367085, 0.5035268961367085, 13.00945052986233, 117.085054768761, 0.1654539549364043, 0.1654539549364043)
$row8 = @("sCs", "Il15", "Il2ra", "ECs", 3, 1, 0.9393656666666668, 2.818097, 0.03599963817648871, 0.03599963817648871, 2, 0.6666666666666666, 0.3213843333333333, 0.964153, 0.1066545648432073, 0.1066545648432073, 0.3018974085378889, 2.717076676841, 0.003839525744226314, 0.003839525744226315)
$row9 = @("sCs", "Il15", "Il2ra", "FAPs", 3, 1, 0.9393656666666668, 2.818097, 0.03599963817648871, 0.03599963817648871, 3, 1, 1.174648, 3.523944, 0.3898185390200842, 0.3898185390200842, 1.103424001618667, 9.930816014568002, 0.01403332635921047, 0.01403332635921048)
$row10 = @("sCs", "Il15", "Il2ra", "sCs", 3, 1, 0.9393656666666668, 2.818097, 0.03599963817648871, 0.03599963817648871, 3, 1, 1.517287666666667, 4.551863, 0.5035268961367085, 0.5035268961367085, 1.425287940523444, 12.827591464711, 0.01812678607305192, 0.01812678607305192)

$data = @($row2, $row3, $row4, $row5, $row6, $row7, $row8, $row9, $row10)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $rowData[$j]
    }
}

